# Update setpoints in USE_TYPE_PROPERTIES / INDOOR_COMFORT sheet
# Increase Cooling setpoint (Tcs_set_C, column B) by 1C
# Decrease Heating setpoint (Ths_set_C, column C) by 1C
# EXCEPT for HOSPITAL and COOLROOM rows, which remain unchanged

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INDOOR_COMFORT")

# Rows to skip (no change): HOSPITAL (row 12), COOLROOM (row 17)
$skipRows = @(12, 17)

for ($row = 2; $row -le 21; $row++) {
    if ($skipRows -contains $row) {
        continue
    }

    $tcs = $ws.Cells.Item($row, 2).Value2   # column B: Tcs_set_C (cooling setpoint)
    $ths = $ws.Cells.Item($row, 3).Value2   # column C: Ths_set_C (heating setpoint)

    $ws.Cells.Item($row, 2).Value = $tcs + 1
    $ws.Cells.Item($row, 3).Value = $ths - 1
}

# Match the author's final UI state (active cell / zoom) on save
$ws.Activate() | Out-Null
$ws.Range("B21").Select() | Out-Null
$excel.ActiveWindow.Zoom = 110

$wb.Save()
